$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the JO No. value that was previously entered in C10 (merged C10:E10)
$ws.Range("C10").Value = ""

# Remove the three extra "Sample scope of work" rows (16:18) that were
# mistakenly left in when the Notes rows were inserted. Deleting these
# entire rows shifts the Notes section (and everything below) up by 3 rows.
$ws.Rows("16:18").Delete()

# Update the view to reflect where the user ended up after the edit.
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("A16:XFD18").Select()
